# TimeManagement.xlsx update:
#  - minimize the workbook window
#  - add three new "Task 5" rows (test-feedback loop tasks), each costing 1 hr

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Best-effort: mark the workbook window as minimized (mirrors the
# `minimized="1"` flag added to <workbookView> in workbook.xml).
try {
    $wb.Windows.Item(1).WindowState = -4140   # xlMinimized
} catch {
    # Window-state isn't modelled everywhere; ignore if unsupported.
}

# New rows appended right after the existing "Task 4: Week 3" block
# (which currently ends at row 24, B24 = "Write up log for testers" / "1 hr").
$ws.Range("A25").Value = "Find up to 10 students to test app"
$ws.Range("B25").Value = "1 hr"

$ws.Range("A26").Value = "Keep a log of suggestions"
$ws.Range("B26").Value = "1 hr"

$ws.Range("A27").Value = "Listen to feedback and make necessary changes"
$ws.Range("B27").Value = "1 hr"

# Mirror the author's last selection landing on the final new cell.
$ws.Range("B27").Select() | Out-Null
